# The template's currency column header changes from "euro" to "sek"
# (commit: "changed templates to sek"). That header lives in cell C1
# of the single worksheet ("Blad1"), which holds the shared-string "euro".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C1").Value = "sek"
